# Apply updated cryptocurrency price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    # Leading apostrophe forces Excel to store the literal as text (matches
    # the workbook's existing inline-string cells) even when it looks numeric;
    # resetting the style back to Normal avoids leaving a quote-prefix format.
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell 'D2' '41.752.95'
Set-TextCell 'E2' '  -1.95%  '

Set-TextCell 'D3' '2.278.84'
Set-TextCell 'E3' '  -2.93%  '

Set-TextCell 'E4' '  +0.00%  '

Set-TextCell 'D5' '315.24'
Set-TextCell 'E5' '  +0.45%  '

Set-TextCell 'D6' '102.39'
Set-TextCell 'E6' '  -5.40%  '

Set-TextCell 'E7' '  -0.96%  '

Set-TextCell 'E8' '  -0.05%  '

Set-TextCell 'E9' '  -2.41%  '

Set-TextCell 'D10' '38.85'
Set-TextCell 'E10' '  -5.58%  '

Set-TextCell 'D11' '0.0903'
Set-TextCell 'E11' '  -2.26%  '

Set-TextCell 'D12' '8.26'
Set-TextCell 'E12' '  -3.43%  '

Set-TextCell 'E13' '  -0.19%  '

Set-TextCell 'D14' '0.959'
Set-TextCell 'E14' '  -3.28%  '

Set-TextCell 'D15' '15.21'
Set-TextCell 'E15' '  -4.54%  '

Set-TextCell 'D16' '2.625.70'
Set-TextCell 'E16' '  -2.80%  '

Set-TextCell 'D17' '2.283.89'
Set-TextCell 'E17' '  -3.05%  '

Set-TextCell 'D18' '41.745.23'
Set-TextCell 'E18' '  -1.87%  '

Set-TextCell 'D19' '7.53'
Set-TextCell 'E19' '  -1.24%  '

Set-TextCell 'E20' '  -0.90%  '

Set-TextCell 'D21' '284.95'

Set-TextCell 'D22' '73.52'
Set-TextCell 'E22' '  -3.18%  '

Set-TextCell 'D23' '3.54'
Set-TextCell 'E23' '  -1.72%  '

Set-TextCell 'E24' '  -2.08%  '

Set-TextCell 'D25' '9.90'
Set-TextCell 'E25' '  +5.56%  '

Set-TextCell 'E26' '  +0.72%  '

Set-TextCell 'D27' '10.72'
Set-TextCell 'E27' '  -5.65%  '

Set-TextCell 'D28' '2.29'
Set-TextCell 'E28' '  +3.44%  '

Set-TextCell 'D29' '22.97'
Set-TextCell 'E29' '  +1.09%  '

Set-TextCell 'D30' '162.85'
Set-TextCell 'E30' '  -5.79%  '

Set-TextCell 'D31' '34.56'
Set-TextCell 'E31' '  -5.48%  '

Set-TextCell 'D32' '0.0873'
Set-TextCell 'E32' '  -1.79%  '

Set-TextCell 'D33' '2.92'
Set-TextCell 'E33' '  +1.69%  '

Set-TextCell 'D34' '5.80'
Set-TextCell 'E34' '  -4.19%  '

Set-TextCell 'E35' '  -0.24%  '

Set-TextCell 'D36' '0.116'
Set-TextCell 'E36' '  -7.31%  '

Set-TextCell 'E37' '  -1.44%  '

Set-TextCell 'E38' '  +8.72%  '

Set-TextCell 'E39' '  -4.25%  '

Set-TextCell 'D40' '3.59'
Set-TextCell 'E40' '  -8.71%  '

Set-TextCell 'D41' '103.34'
Set-TextCell 'E41' '  +22.38%  '

Set-TextCell 'D42' '1.45'
Set-TextCell 'E42' '  -0.95%  '

Set-TextCell 'D43' '69.33'
Set-TextCell 'E43' '  -1.95%  '

Set-TextCell 'E44' '  +0.28%  '

Set-TextCell 'D45' '0.224'
Set-TextCell 'E45' '  -5.10%  '

Set-TextCell 'D46' '114.92'
Set-TextCell 'E46' '  +3.15%  '

Set-TextCell 'D47' '11.93'
Set-TextCell 'E47' '  +0.03%  '

Set-TextCell 'D48' '8.97'
Set-TextCell 'E48' '  -1.95%  '

Set-TextCell 'D49' '76.10'
Set-TextCell 'E49' '  +1.49%  '

Set-TextCell 'E50' '  -3.00%  '

Set-TextCell 'E51' '  -1.78%  '
